$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record was inserted right before the existing row 196,
# pushing every following row (196-246 -> 197-247) down by one. Replicate
# that with a native row insert, then populate the newly opened row 196
# with the record's data.
$ws.Rows.Item(196).Insert()

$ws.Range("A196").Value = 8
$ws.Range("B196").Value = "Terminal La Palmera de La Serena"
$ws.Range("C196").Value = "Coquimbo"
$ws.Range("D196").Value = 44841
$ws.Range("E196").Value = 4
$ws.Range("F196").Value = 100112037
$ws.Range("G196").Value = "Cebollín"
$ws.Range("H196").Value = "Sin especificar"
$ws.Range("I196").Value = "Primera"
$ws.Range("J196").Value = 1800
$ws.Range("K196").Value = 1400
$ws.Range("L196").Value = 1600
$ws.Range("M196").Value = 1500
$ws.Range("N196").Value = "`$/paquete 6 unidades"
$ws.Range("O196").Value = "Provincia del Elquí"
$ws.Range("P196").Value = 250
$ws.Range("Q196").Value = 6
$ws.Range("R196").Value = "Hortaliza"
